$d = $word.ActiveDocument
$d.Content.Find.Execute("Medicaid Services", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Medicaid Services", 2)
Write-Output "done"
